# Insert a new data row at row 18. Excel will automatically push the
# existing rows 18..80 down to 19..81 (copying formatting as appropriate),
# which matches the "each row shifts down by one" pattern seen in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 45125
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = 100112043
$ws.Cells.Item(18, 7).Value = "Pepino dulce"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 150
$ws.Cells.Item(18, 11).Value = 16000
$ws.Cells.Item(18, 12).Value = 17000
$ws.Cells.Item(18, 13).Value = 16400
$ws.Cells.Item(18, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 911
$ws.Cells.Item(18, 17).Value = 18
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date/time style used by the rest of
# column D (style index 2 = numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
